$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.444.37'
$ws.Range("D3").Value = '1.828.61'
$ws.Range("E3").Value = '  -2.24%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -1.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '331.04'
$ws.Range("E5").Value = '  -1.57%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -1.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4582'
$ws.Range("E7").Value = '  -2.31%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3824'
$ws.Range("E8").Value = '  -2.65%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.46'
$ws.Range("E9").Value = '  +3.79%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07892'
$ws.Range("E10").Value = '  -1.51%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9670'
$ws.Range("E11").Value = '  -3.90%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.09'
$ws.Range("E12").Value = '  -3.48%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.875'
$ws.Range("E13").Value = '  -2.10%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.825.31'
$ws.Range("E14").Value = '  -2.75%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.091'
$ws.Range("E15").Value = '  -2.40%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.002'
$ws.Range("E16").Value = '  -1.14%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '89.70'
$ws.Range("E17").Value = '  +0.99%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06602'
$ws.Range("E18").Value = '  -2.35%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.00001025'
$ws.Range("E19").Value = '  -2.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.15'
$ws.Range("E20").Value = '  -0.76%  '
$ws.Range("D22").Value = '27.431.06'
$ws.Range("E22").Value = '  -1.16%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.329'
$ws.Range("E23").Value = '  -2.92%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.82'
$ws.Range("E24").Value = '  -1.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.282'
$ws.Range("E25").Value = '  -1.70%  '
$ws.Range("D26").Value = '2.048.12'
$ws.Range("E26").Value = '  -2.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '156.15'
$ws.Range("E27").Value = '  -1.99%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.41'
$ws.Range("E28").Value = '  -2.12%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.069'
$ws.Range("E29").Value = '  -3.77%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.295'
$ws.Range("E30").Value = '  -3.03%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '118.49'
$ws.Range("E31").Value = '  -2.88%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09310'
$ws.Range("E32").Value = '  -2.25%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9373'
$ws.Range("E33").Value = '  -4.66%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.574'
$ws.Range("E34").Value = '  -1.68%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.239'
$ws.Range("E35").Value = '  -2.03%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.328'
$ws.Range("E36").Value = '  -0.83%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05922'
$ws.Range("E37").Value = '  -2.46%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02177'
$ws.Range("E38").Value = '  -2.91%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.117'
$ws.Range("E39").Value = '  -2.09%  '
$ws.Range("E40").Value = '  -1.04%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.148'
$ws.Range("E41").Value = '  -3.99%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5776'
$ws.Range("E42").Value = '  -3.61%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1823'
$ws.Range("E43").Value = '  -3.72%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '9.990'
$ws.Range("E44").Value = '  -3.28%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.282'
$ws.Range("E45").Value = '  +2.83%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '11.98'
$ws.Range("E46").Value = '  -1.72%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5433'
$ws.Range("E47").Value = '  -4.29%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.876'
$ws.Range("E48").Value = '  -3.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '110.17'
$ws.Range("E49").Value = '  -1.95%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06582'
$ws.Range("E50").Value = '  -2.72%  '
$ws.Range("E51").Value = '  -33.64%  '
